# Realestate Update resale numbers 2023-06-26 14:42
# Append a new data row (row 78) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 78

# Columns A (Date) and D (Week) look numeric/date-like to Excel's auto-detection,
# so a leading apostrophe keeps them stored as literal text (matching columns B/C
# and the rest of the sheet's text-formatted Date/Time/Weekday/Week columns).
$ws.Cells.Item($row, 1).Value = "'2023-06-26"
$ws.Cells.Item($row, 2).Value = "14:42:09"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "'26"

$ws.Cells.Item($row, 5).Value = 122745
$ws.Cells.Item($row, 6).Value = 134381
$ws.Cells.Item($row, 7).Value = 163465
$ws.Cells.Item($row, 8).Value = 133635
$ws.Cells.Item($row, 9).Value = 177312
$ws.Cells.Item($row, 10).Value = 115068
$ws.Cells.Item($row, 11).Value = 203471
$ws.Cells.Item($row, 12).Value = 226046
$ws.Cells.Item($row, 13).Value = 176244
$ws.Cells.Item($row, 14).Value = 104319
$ws.Cells.Item($row, 15).Value = 39613
$ws.Cells.Item($row, 16).Value = 33789
$ws.Cells.Item($row, 17).Value = 52185
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35888
$ws.Cells.Item($row, 20).Value = -1
